$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 5 and 6 as completed by setting column B to 1 (same as other
# finished tasks above them).
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1

# Highlight A7 with the same "done" fill/style used by A1:A3, A8 (fillId=2,
# solid green 92D050) so it picks up the existing style index (s="1").
$ws.Range("A7").Interior.Color = $ws.Range("A1").Interior.Color

# Move the active selection from A10 to A8.
$ws.Range("A8").Select()
